$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits right
#    after "Towers" in the title paragraph.
# ---------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------
# 2) "Orientador:" -> "Supervisor" + ":"  (two runs, same formatting)
#    Paragraph is the one whose text starts with "Orientador:".
# ---------------------------------------------------------------
$pSup = $d.Paragraphs(20)
$start = $pSup.Range.Start
$wordRange = $d.Range($start, $start + 10)   # "Orientador"
$wordRange.Text = "Supervisor"
# Force a run split right after the replaced word (so the trailing
# ":" stays in its own run) by dropping a temporary bookmark there
# and then removing it again - the split survives the removal.
$splitPos = $start + 10
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("TmpSplit1", $splitRange) | Out-Null
$d.Bookmarks("TmpSplit1").Delete()

# ---------------------------------------------------------------
# 3) "Coorientador:" -> "Co-supervisor:" with the _GoBack bookmark
#    re-inserted between "-s" and "upervisor".
# ---------------------------------------------------------------
$pCoSup = $d.Paragraphs(23)
$coStart = $pCoSup.Range.Start
# "Co" (2 chars) stays as is.
# "o" (1 char) -> "-s"
$oRange = $d.Range($coStart + 2, $coStart + 3)
$oRange.Text = "-s"
# "rientador:" (10 chars) -> "upervisor" + ":"
$restStart = $coStart + 2 + 2   # length of "-s" is 2
$restRange = $d.Range($restStart, $restStart + 10)
$restRange.Text = "upervisorX"
# Re-insert the _GoBack bookmark right before "upervisor".
$bmPos = $restStart
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
# Force the trailing ":" into its own run the same way as above,
# then fix the placeholder "X" back into the final ":".
$tailPos = $restStart + 9
$tailSplit = $d.Range($tailPos, $tailPos)
$d.Bookmarks.Add("TmpSplit2", $tailSplit) | Out-Null
$d.Bookmarks("TmpSplit2").Delete()
$xRange = $d.Range($tailPos, $tailPos + 1)
$xRange.Text = ":"

# ---------------------------------------------------------------
# 4) Drop the highlighted placeholder note
#    " (Mes e ano em Arial 11pt maiusculas pequenas na ultima linha)"
#    leaving its paragraph empty.
# ---------------------------------------------------------------
$pNote = $d.Paragraphs(27)
$noteRange = $pNote.Range
$noteRange.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
$noteRange.Text = ""

# ---------------------------------------------------------------
# 5) Year fix: "junho de 2015" -> "junho de 2017"
# ---------------------------------------------------------------
$pDate = $d.Paragraphs(28)
$dateStart = $pDate.Range.Start
$digitRange = $d.Range($dateStart + 12, $dateStart + 13)
$digitRange.Text = "7"
